# Auto-generated edit script applying scheduled-runner value updates
# to the per-job Excalibur_Profits leve-profit tables across all 8 sheets.
$wb = $excel.ActiveWorkbook

# ALC row 40
$ws = $wb.Worksheets.Item(1)
$ws.Range("H40").Value = 4817.6665
$ws.Range("J40").Value = 4626.5
$ws.Range("L40").Value = 4626.5
$ws.Range("N40").Value = -4976.5

# ALC row 43
$ws = $wb.Worksheets.Item(1)
$ws.Range("H43").Value = 2449.9
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 2562.375
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 2562.375
$ws.Range("M43").Value = -1931
$ws.Range("N43").Value = -2700.375

# ALC row 113
$ws = $wb.Worksheets.Item(1)
$ws.Range("H113").Value = 4950
$ws.Range("I113").Value = 4950
$ws.Range("K113").Value = 4950
$ws.Range("M113").Value = -1696

# ALC row 136
$ws = $wb.Worksheets.Item(1)
$ws.Range("H136").Value = 78948.75
$ws.Range("J136").Value = 78948.75
$ws.Range("L136").Value = 78948.75
$ws.Range("N136").Value = -89148.75

# ALC row 138
$ws = $wb.Worksheets.Item(1)
$ws.Range("H138").Value = 3247.6272
$ws.Range("J138").Value = 4402.5947
$ws.Range("L138").Value = 13207.7841
$ws.Range("N138").Value = -23487.7841

# ARM row 11
$ws = $wb.Worksheets.Item(2)
$ws.Range("H11").Value = 1001700
$ws.Range("J11").Value = 3400
$ws.Range("L11").Value = 3400
$ws.Range("N11").Value = -3688

# ARM row 24
$ws = $wb.Worksheets.Item(2)
$ws.Range("H24").Value = 29677.5
$ws.Range("J24").Value = 29677.5
$ws.Range("L24").Value = 29677.5
$ws.Range("N24").Value = -30425.5

# ARM row 32
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 8774403
$ws.Range("I32").Value = 10639712
$ws.Range("K32").Value = 10639712
$ws.Range("M32").Value = -10639425

# ARM row 74
$ws = $wb.Worksheets.Item(2)
$ws.Range("H74").Value = 2474.5173
$ws.Range("I74").Value = 1334.8125
$ws.Range("J74").Value = 3877.2307
$ws.Range("K74").Value = 1334.8125
$ws.Range("L74").Value = 3877.2307
$ws.Range("M74").Value = -460.8125
$ws.Range("N74").Value = -5625.2307

# ARM row 77
$ws = $wb.Worksheets.Item(2)
$ws.Range("H77").Value = 2474.5173
$ws.Range("I77").Value = 1334.8125
$ws.Range("J77").Value = 3877.2307
$ws.Range("K77").Value = 6674.0625
$ws.Range("L77").Value = 19386.1535
$ws.Range("M77").Value = -2306.0625
$ws.Range("N77").Value = -28122.1535

# ARM row 100
$ws = $wb.Worksheets.Item(2)
$ws.Range("H100").Value = 29677.5
$ws.Range("J100").Value = 29677.5
$ws.Range("L100").Value = 29677.5
$ws.Range("N100").Value = -31841.5

# ARM row 110
$ws = $wb.Worksheets.Item(2)
$ws.Range("H110").Value = 1006.1667
$ws.Range("I110").Value = 831.7143
$ws.Range("J110").Value = 1616.75
$ws.Range("K110").Value = 831.7143
$ws.Range("L110").Value = 1616.75
$ws.Range("M110").Value = 1213.2857
$ws.Range("N110").Value = -5706.75

# ARM row 122
$ws = $wb.Worksheets.Item(2)
$ws.Range("H122").Value = 1332.1666
$ws.Range("I122").Value = 1000.1
$ws.Range("K122").Value = 3000.3
$ws.Range("M122").Value = -550.3000000000002

# ARM row 134
$ws = $wb.Worksheets.Item(2)
$ws.Range("H134").Value = 68850
$ws.Range("J134").Value = 68850
$ws.Range("L134").Value = 68850
$ws.Range("N134").Value = -78990

# ARM row 139
$ws = $wb.Worksheets.Item(2)
$ws.Range("H139").Value = 174499.5
$ws.Range("J139").Value = 174499.5
$ws.Range("L139").Value = 174499.5
$ws.Range("N139").Value = -184779.5

# BSM row 5
$ws = $wb.Worksheets.Item(3)
$ws.Range("H5").Value = 4303.778
$ws.Range("I5").Value = 374.83334
$ws.Range("J5").Value = 12161.667
$ws.Range("K5").Value = 374.83334
$ws.Range("L5").Value = 12161.667
$ws.Range("M5").Value = -261.83334
$ws.Range("N5").Value = -12387.667

# BSM row 86
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 3556.8572
$ws.Range("I86").Value = 3500
$ws.Range("K86").Value = 3500
$ws.Range("M86").Value = -2377

# BSM row 89
$ws = $wb.Worksheets.Item(3)
$ws.Range("H89").Value = 3556.8572
$ws.Range("I89").Value = 3500
$ws.Range("K89").Value = 17500
$ws.Range("M89").Value = -11884

# BSM row 100
$ws = $wb.Worksheets.Item(3)
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = ""

# CRP row 16
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 1746.1666
$ws.Range("I16").Value = 1584.6
$ws.Range("J16").Value = 1861.5714
$ws.Range("K16").Value = 1584.6
$ws.Range("L16").Value = 1861.5714
$ws.Range("M16").Value = -1297.6
$ws.Range("N16").Value = -2435.5714

# CRP row 18
$ws = $wb.Worksheets.Item(4)
$ws.Range("H18").Value = 99994
$ws.Range("J18").Value = 99994
$ws.Range("L18").Value = 99994
$ws.Range("N18").Value = -100454

# CRP row 22
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 1899
$ws.Range("J22").Value = 1899
$ws.Range("L22").Value = 1899
$ws.Range("N22").Value = -2599

# CRP row 31
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 67098
$ws.Range("I31").Value = 86054.71000000001
$ws.Range("J31").Value = 7068.4165
$ws.Range("K31").Value = 86054.71000000001
$ws.Range("L31").Value = 7068.4165
$ws.Range("M31").Value = -85759.71000000001
$ws.Range("N31").Value = -7658.4165

# CRP row 34
$ws = $wb.Worksheets.Item(4)
$ws.Range("H34").Value = 67098
$ws.Range("I34").Value = 86054.71000000001
$ws.Range("J34").Value = 7068.4165
$ws.Range("K34").Value = 86054.71000000001
$ws.Range("L34").Value = 7068.4165
$ws.Range("M34").Value = -85852.71000000001
$ws.Range("N34").Value = -7472.4165

# CRP row 41
$ws = $wb.Worksheets.Item(4)
$ws.Range("H41").Value = 11999
$ws.Range("I41").Value = 1000
$ws.Range("J41").Value = 33997
$ws.Range("K41").Value = 1000
$ws.Range("L41").Value = 33997
$ws.Range("M41").Value = -572
$ws.Range("N41").Value = -34853

# CRP row 50
$ws = $wb.Worksheets.Item(4)
$ws.Range("H50").Value = 39500
$ws.Range("J50").Value = 46666.668
$ws.Range("L50").Value = 46666.668
$ws.Range("N50").Value = -47916.668

# CRP row 62
$ws = $wb.Worksheets.Item(4)
$ws.Range("H62").Value = 1825
$ws.Range("I62").Value = 1900
$ws.Range("J62").Value = 1750
$ws.Range("K62").Value = 1900
$ws.Range("L62").Value = 1750
$ws.Range("M62").Value = -1276
$ws.Range("N62").Value = -2998

# CRP row 65
$ws = $wb.Worksheets.Item(4)
$ws.Range("H65").Value = 1825
$ws.Range("I65").Value = 1900
$ws.Range("J65").Value = 1750
$ws.Range("K65").Value = 9500
$ws.Range("L65").Value = 8750
$ws.Range("M65").Value = -6380
$ws.Range("N65").Value = -14990

# CRP row 68
$ws = $wb.Worksheets.Item(4)
$ws.Range("H68").Value = 78897.25
$ws.Range("J68").Value = 78897.25
$ws.Range("L68").Value = 78897.25
$ws.Range("N68").Value = -80395.25

# CRP row 71
$ws = $wb.Worksheets.Item(4)
$ws.Range("H71").Value = 78897.25
$ws.Range("J71").Value = 78897.25
$ws.Range("L71").Value = 236691.75
$ws.Range("N71").Value = -244179.75

# CRP row 113
$ws = $wb.Worksheets.Item(4)
$ws.Range("H113").Value = 1746.1666
$ws.Range("I113").Value = 1584.6
$ws.Range("J113").Value = 1861.5714
$ws.Range("K113").Value = 1584.6
$ws.Range("L113").Value = 1861.5714
$ws.Range("M113").Value = 585.4000000000001
$ws.Range("N113").Value = -6201.5714

# CRP row 122
$ws = $wb.Worksheets.Item(4)
$ws.Range("H122").Value = 1780.2858
$ws.Range("I122").Value = 835.1667
$ws.Range("J122").Value = 2489.125
$ws.Range("K122").Value = 2505.5001
$ws.Range("L122").Value = 7467.375
$ws.Range("M122").Value = -55.5001000000002
$ws.Range("N122").Value = -12367.375

# CRP row 132
$ws = $wb.Worksheets.Item(4)
$ws.Range("H132").Value = 43104600
$ws.Range("I132").Value = 37038160
$ws.Range("K132").Value = 111114480
$ws.Range("M132").Value = -111111950

# GSM row 14
$ws = $wb.Worksheets.Item(6)
$ws.Range("H14").Value = 8995011
$ws.Range("I14").Value = 8666682
$ws.Range("K14").Value = 8666682
$ws.Range("M14").Value = -8666514

# GSM row 17
$ws = $wb.Worksheets.Item(6)
$ws.Range("H17").Value = 4750
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 7500
$ws.Range("K17").Value = 2000
$ws.Range("L17").Value = 7500
$ws.Range("M17").Value = -1832
$ws.Range("N17").Value = -7836

# GSM row 102
$ws = $wb.Worksheets.Item(6)
$ws.Range("H102").Value = 3278.484
$ws.Range("I102").Value = 3134.5557
$ws.Range("K102").Value = 3134.5557
$ws.Range("M102").Value = -1512.5557

# GSM row 113
$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 2001.2
$ws.Range("I113").Value = 1669
$ws.Range("J113").Value = 2499.5
$ws.Range("K113").Value = 1669
$ws.Range("L113").Value = 2499.5
$ws.Range("M113").Value = 501
$ws.Range("N113").Value = -6839.5

# GSM row 139
$ws = $wb.Worksheets.Item(6)
$ws.Range("H139").Value = 102305.75
$ws.Range("J139").Value = 102305.75
$ws.Range("L139").Value = 102305.75
$ws.Range("N139").Value = -112585.75

# GSM row 140
$ws = $wb.Worksheets.Item(6)
$ws.Range("H140").Value = 97975
$ws.Range("J140").Value = 97975
$ws.Range("L140").Value = 97975
$ws.Range("N140").Value = -108335

# LTW row 12
$ws = $wb.Worksheets.Item(7)
$ws.Range("H12").Value = 368.6
$ws.Range("I12").Value = 368.6
$ws.Range("K12").Value = 368.6
$ws.Range("M12").Value = -198.6

# LTW row 46
$ws = $wb.Worksheets.Item(7)
$ws.Range("H46").Value = 18069.834
$ws.Range("I46").Value = 38209.5
$ws.Range("K46").Value = 38209.5
$ws.Range("M46").Value = -38021.5

# LTW row 61
$ws = $wb.Worksheets.Item(7)
$ws.Range("H61").Value = 2479.0857
$ws.Range("I61").Value = 1499.64
$ws.Range("J61").Value = 4927.7
$ws.Range("K61").Value = 1499.64
$ws.Range("L61").Value = 4927.7
$ws.Range("M61").Value = -1297.64
$ws.Range("N61").Value = -5331.7

# LTW row 68
$ws = $wb.Worksheets.Item(7)
$ws.Range("H68").Value = 4283.75
$ws.Range("I68").Value = 2599.8
$ws.Range("K68").Value = 2599.8
$ws.Range("M68").Value = -1850.8

# LTW row 71
$ws = $wb.Worksheets.Item(7)
$ws.Range("H71").Value = 4283.75
$ws.Range("I71").Value = 2599.8
$ws.Range("K71").Value = 12999
$ws.Range("M71").Value = -9255

# LTW row 113
$ws = $wb.Worksheets.Item(7)
$ws.Range("H113").Value = 2479.0857
$ws.Range("I113").Value = 1499.64
$ws.Range("J113").Value = 4927.7
$ws.Range("K113").Value = 1499.64
$ws.Range("L113").Value = 4927.7
$ws.Range("M113").Value = 670.3599999999999
$ws.Range("N113").Value = -9267.700000000001

# LTW row 122
$ws = $wb.Worksheets.Item(7)
$ws.Range("H122").Value = 40304.73
$ws.Range("I122").Value = 1489.8096
$ws.Range("J122").Value = 203327.4
$ws.Range("K122").Value = 4469.4288
$ws.Range("L122").Value = 609982.2
$ws.Range("M122").Value = -2019.4288
$ws.Range("N122").Value = -614882.2

# WVR row 132
$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 4764497
$ws.Range("I132").Value = 6898770
$ws.Range("K132").Value = 20696310
$ws.Range("M132").Value = -20693780
